# Business logic to plugins.
# - MQTT processing to iotgate.
# - SoC temperature measurement to system.
#
# Applies the semantic re-classification of the "soc" category (now "system")
# and splits the combined "on"/"off" + "perc"/"temp" rows into explicit
# percon/percoff/tempon/tempoff parameter rows on the "Semanticke clenenie" sheet.
# Also adds a new "Comment" column and renames the value column header.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- fan / state: split perc + on/off into percon/percoff ---------------
$ws.Range("C5").Value = "percon"
$ws.Range("D5").ClearContents()
$ws.Range("E5").Value = 95

$ws.Range("C6").Value = "percoff"
$ws.Range("D6").ClearContents()
$ws.Range("E6").Value = 60

# --- fan / state: split temp + on/off into tempon/tempoff ---------------
$ws.Range("C7").Value = "tempon"
$ws.Range("D7").ClearContents()

$ws.Range("C8").Value = "tempoff"
$ws.Range("D8").ClearContents()

# --- Header row ---------------------------------------------------------
$ws.Range("F1").Value = "Comment"
$ws.Range("E1").Value = "payload (value)"

# --- soc -> system (rows 2 and 3) --------------------------------------
$ws.Range("A2").Value = "system"
$ws.Range("A3").Value = "system"

# --- fan / cmd: split perc + on/off into percon/percoff -----------------
$ws.Range("C10").Value = "percon"
$ws.Range("D10").ClearContents()

$ws.Range("C11").Value = "percoff"
$ws.Range("D11").ClearContents()

# --- fan / cmd: split temp + on/off into tempon/tempoff -----------------
$ws.Range("C12").Value = "tempon"
$ws.Range("D12").ClearContents()

$ws.Range("C13").Value = "tempoff"
$ws.Range("D13").ClearContents()

# --- selection moves to D3 ------------------------------------------------
$ws.Range("D3").Select()
